# Apply the metadata refresh edit described by the diff:
#  - Update the "Date" value on the Metadata sheet
#  - Update the "System URI" values on the two Include sheets to point at
#    the new interop.esante.gouv.fr terminology host

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value (row 8, column B) ---
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# --- Include #0 sheet: update System URI value (row 4, column B) ---
$include0 = $wb.Worksheets.Item("Include #0")
$include0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R20-Pays"

# --- Include #1 sheet: update System URI value (row 4, column B) ---
$include1 = $wb.Worksheets.Item("Include #1")
$include1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R13-CommuneOM"
